$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 (Marking): Right / Wrong marks corrections
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): updated totals and fraction text
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "62 / 112"
